$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Pre-intern "SOMME DES LIGNES" ahead of the "destinataire client" fix so
#    the shared-strings table ends up ordered the same way as the target
#    workbook (cosmetic only - cell references stay correct either way).
# ---------------------------------------------------------------------------
$ws.Range("ZZ100").Value = "SOMME DES LIGNES"
$ws.Range("ZZ100").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 1. Fix the typo in the "destinataire client" column header (D4).
#    ("destinatire client" -> "destinataire client")
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "destinataire client"

# ---------------------------------------------------------------------------
# 2. Title bar (B2:G2, merged) gets a dashed box border around it.
#    Left edge (B2): left + top + bottom mediumDashed
#    Middle (C2:F2): top + bottom mediumDashed
#    Right edge (G2): right + top + bottom mediumDashed
# ---------------------------------------------------------------------------
$b2 = $ws.Range("B2")
$b2.Borders.Item(7).LineStyle = -4115
$b2.Borders.Item(7).Weight = -4138
$b2.Borders.Item(8).LineStyle = -4115
$b2.Borders.Item(8).Weight = -4138
$b2.Borders.Item(9).LineStyle = -4115
$b2.Borders.Item(9).Weight = -4138

$mid2 = $ws.Range("C2:F2")
$mid2.Borders.Item(8).LineStyle = -4115
$mid2.Borders.Item(8).Weight = -4138
$mid2.Borders.Item(9).LineStyle = -4115
$mid2.Borders.Item(9).Weight = -4138

$g2 = $ws.Range("G2")
$g2.Borders.Item(10).LineStyle = -4115
$g2.Borders.Item(10).Weight = -4138
$g2.Borders.Item(8).LineStyle = -4115
$g2.Borders.Item(8).Weight = -4138
$g2.Borders.Item(9).LineStyle = -4115
$g2.Borders.Item(9).Weight = -4138

$ws.Rows.Item(2).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3. Row 3 height only changes (formats already match).
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 24

# ---------------------------------------------------------------------------
# 4. Table header row (B4:G4): double top / thick bottom border,
#    horizontally centered, taller row.
# ---------------------------------------------------------------------------
$hdr = $ws.Range("B4:G4")
$hdr.Borders.Item(8).LineStyle = -4119
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = 4
$hdr.HorizontalAlignment = -4108
$ws.Rows.Item(4).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 5. Row 5 (data row) only gets a taller height (formats already match).
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 6. New summation row 6.
#    B6:D6 -> thick top border, vertically centered, blank
#    E6:F6 -> merged, thick top border, 14pt font, centered, "SOMME DES LIGNES"
#    G6    -> thick top border, vertically centered, currency format, =SUM(G5:G5)
# ---------------------------------------------------------------------------
$left6 = $ws.Range("B6:D6")
$left6.Borders.Item(8).LineStyle = 1
$left6.Borders.Item(8).Weight = 4
$left6.VerticalAlignment = -4108

$ws.Range("E6:F6").Merge() | Out-Null
$lbl6 = $ws.Range("E6:F6")
$lbl6.Borders.Item(8).LineStyle = 1
$lbl6.Borders.Item(8).Weight = 4
$lbl6.Font.Size = 14
$lbl6.HorizontalAlignment = -4108
$lbl6.VerticalAlignment = -4108
$ws.Range("E6").Value = "SOMME DES LIGNES"

$g6 = $ws.Range("G6")
$g6.Borders.Item(8).LineStyle = 1
$g6.Borders.Item(8).Weight = 4
$g6.VerticalAlignment = -4108
$g6.NumberFormat = "#,##0.00\ ""€"""
$g6.Formula = "=SUM(G5:G5)"

$ws.Rows.Item(6).RowHeight = 19.5

# ---------------------------------------------------------------------------
# 7. Selection moves to G13.
# ---------------------------------------------------------------------------
$ws.Range("G13").Select() | Out-Null
